$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.133.77'
$ws.Range('E2').Value = '  +1.40%  '
$ws.Range('D3').Value = '4.023.66'
$ws.Range('E3').Value = '  +0.68%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '535.84'
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.52'
$ws.Range('E6').Value = '  +3.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.707'
$ws.Range('E7').Value = '  +13.49%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.753'
$ws.Range('E9').Value = '  +1.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.172'
$ws.Range('E10').Value = '  -3.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000324'
$ws.Range('E11').Value = '  -7.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.52'
$ws.Range('E12').Value = '  +10.79%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.71'
$ws.Range('E13').Value = '  +1.74%  '
$ws.Range('D14').Value = '4.667.60'
$ws.Range('E14').Value = '  +1.36%  '
$ws.Range('D15').Value = '4.040.98'
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.08'
$ws.Range('E16').Value = '  -1.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '20.52'
$ws.Range('E17').Value = '  -3.92%  '
$ws.Range('E18').Value = '  -0.70%  '
$ws.Range('E19').Value = '  -2.70%  '
$ws.Range('D20').Value = '72.036.22'
$ws.Range('E20').Value = '  +1.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '428.58'
$ws.Range('E21').Value = '  -2.99%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '99.22'
$ws.Range('E22').Value = '  +9.77%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.53'
$ws.Range('E23').Value = '  -0.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.25'
$ws.Range('E24').Value = '  +4.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '14.44'
$ws.Range('E25').Value = '  +1.43%  '
$ws.Range('E26').Value = '  -6.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.86'
$ws.Range('E27').Value = '  +0.84%  '
$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.85'
$ws.Range('E28').Value = '  +2.02%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.95'
$ws.Range('E29').Value = '  -0.45%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.61'
$ws.Range('E30').Value = '  +25.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.46'
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.130'
$ws.Range('E32').Value = '  +2.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '681.31'
$ws.Range('E33').Value = '  -2.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.94'
$ws.Range('E34').Value = '  +2.26%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '66.34'
$ws.Range('E35').Value = '  -1.66%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '42.80'
$ws.Range('E36').Value = '  +6.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.425'
$ws.Range('E37').Value = '  -3.57%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.153'
$ws.Range('E38').Value = '  +2.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.46'
$ws.Range('E39').Value = '  +8.92%  '
$ws.Range('D40').Value = '0.0₃0819'
$ws.Range('E40').Value = '  -12.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.39'
$ws.Range('E41').Value = '  -1.02%  '
$ws.Range('E42').Value = '  -0.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.999'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0487'
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('E45').Value = '  +6.21%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.42'
$ws.Range('E46').Value = '  -2.85%  '
$ws.Range('B47').Value = 'THORChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.64'
$ws.Range('E47').Value = '  +4.23%  '
$ws.Range('E48').Value = '  -10.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.00'
$ws.Range('E49').Value = '  -8.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.36'
$ws.Range('E50').Value = '  -0.58%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '146.29'
$ws.Range('E51').Value = '  +1.43%  '
